$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.439664333333334
$ws.Range("H2").Value = 19.318993
$ws.Range("I2").Value = 0.8861484385177248
$ws.Range("J2").Value = 0.8861484385177248
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 753.7276977616864
$ws.Range("R2").Value = 6783.549279855178
$ws.Range("S2").Value = 0.2875875699114926
$ws.Range("T2").Value = 0.2875875699114926

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.439664333333334
$ws.Range("H3").Value = 19.318993
$ws.Range("I3").Value = 0.8861484385177248
$ws.Range("J3").Value = 0.8861484385177248
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 654.1413433941352
$ws.Range("R3").Value = 5887.272090547217
$ws.Range("S3").Value = 0.2495900308347687
$ws.Range("T3").Value = 0.2495900308347687

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.439664333333334
$ws.Range("H4").Value = 19.318993
$ws.Range("I4").Value = 0.8861484385177248
$ws.Range("J4").Value = 0.8861484385177248
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 914.6048496477139
$ws.Range("R4").Value = 8231.443646829426
$ws.Range("S4").Value = 0.3489708377714635
$ws.Range("T4").Value = 0.3489708377714635

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.8273623333333333
$ws.Range("H5").Value = 2.482087
$ws.Range("I5").Value = 0.1138515614822752
$ws.Range("J5").Value = 0.1138515614822752
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 96.83826274766032
$ws.Range("R5").Value = 871.544364728943
$ws.Range("S5").Value = 0.03694899463128885
$ws.Range("T5").Value = 0.03694899463128885

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.8273623333333333
$ws.Range("H6").Value = 2.482087
$ws.Range("I6").Value = 0.1138515614822752
$ws.Range("J6").Value = 0.1138515614822752
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 84.04349670819377
$ws.Range("R6").Value = 756.391470373744
$ws.Range("S6").Value = 0.03206710468110727
$ws.Range("T6").Value = 0.03206710468110726

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.8273623333333333
$ws.Range("H7").Value = 2.482087
$ws.Range("I7").Value = 0.1138515614822752
$ws.Range("J7").Value = 0.1138515614822752
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 117.5076158186684
$ws.Range("R7").Value = 1057.568542368016
$ws.Range("S7").Value = 0.04483546216987906
$ws.Range("T7").Value = 0.04483546216987905
